# Fixed #366 User content is lost after two generation without edition.
#
# Replace the two <w:fldSimple> fields ("m:usercontent zone1" and
# "m:endusercontent") inside the table cell with the equivalent "complex"
# field representation (separate begin/instrText/separate/end runs), so
# the field result run can carry its own content without losing the
# simple-field wrapper across successive generations.

$d = $word.ActiveDocument

function Get-ParagraphAt($pos) {
    foreach ($p in $d.Paragraphs) {
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

function Replace-SimpleFieldWithComplex($field, $paraAttrsXml, $paraPrXml, $instr) {
    $para = Get-ParagraphAt($field.Code.Start)
    $start = $para.Range.Start
    $end = $para.Range.End

    $bodyXml = "<w:p $paraAttrsXml>$paraPrXml<w:r><w:fldChar w:fldCharType=`"begin`"/></w:r><w:r><w:instrText>$instr</w:instrText></w:r><w:r><w:fldChar w:fldCharType=`"separate`"/></w:r><w:r><w:fldChar w:fldCharType=`"end`"/></w:r></w:p>"

    $packageXml = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $bodyXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r = $d.Range($start, $end)
    $r.InsertXML($packageXml)
}

# Process the later field first so the earlier field's positions stay valid.
$fEnd = $d.Fields.Item(2)
Replace-SimpleFieldWithComplex $fEnd `
    'w:rsidR="00833091" w:rsidRPr="00122CE2" w:rsidRDefault="00122CE2" w:rsidP="00F5495F"' `
    '<w:pPr><w:rPr><w:rStyle w:val="lev"/><w:color w:val="00B050"/></w:rPr></w:pPr>' `
    'm:endusercontent'

$fBegin = $d.Fields.Item(1)
Replace-SimpleFieldWithComplex $fBegin `
    'w:rsidP="004A37BA" w:rsidR="004A37BA" w:rsidRDefault="004A37BA"' `
    '' `
    'm:usercontent zone1'
